$p = $ppt.ActivePresentation

# The deck currently carries the "Integral" theme (ppt/theme/theme2.xml,
# used by the slide master) and an "Office Theme" (ppt/theme/theme1.xml,
# used only by the notes master). The commit swaps the two themes' content,
# which - since the font/format schemes of the two theme parts are already
# identical - amounts to applying the "Office Theme" colour scheme to the
# presentation's (only) design/theme.

$master = $p.SlideMaster
$theme = $master.Theme
$scheme = $theme.ThemeColorScheme

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$scheme.Colors(1).RGB  = 0         # dk1      000000
$scheme.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388   # dk2      44546A
$scheme.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501   # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407     # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308  # accent5  4472C4
$scheme.Colors(10).RGB = 4697456   # accent6  70AD47
$scheme.Colors(11).RGB = 12673797  # hlink    0563C1
$scheme.Colors(12).RGB = 7491477   # folHlink 954F72
